# Bills-RUBIO-Present.xlsx update:
#  - Refresh the Status (col F) values for each bill row to reflect the
#    latest legislative action (Senate/Assembly -> "to Senate"/"to Assembly",
#    committee holds -> "HELD <chamber> - <committee>").
#  - Rename the sheet to the new export timestamp.
#  - Move the active selection to F52 (where the export script left it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Status (column F) values, keyed by row number.
$statusByRow = @{
    2  = "to Senate"
    3  = "Senate - Health"
    4  = "to Senate"
    5  = "to Senate"
    6  = "to Senate"
    7  = "to Senate"
    8  = "HELD Assembly - Education"
    9  = "HELD Senate - Appropriations"
    10 = "to Assembly"
    11 = "to Assembly"
    12 = "to Assembly"
    13 = "to Assembly"
    14 = "to Assembly"
    15 = "to Assembly"
    16 = "to Assembly"
    17 = "to Assembly"
    18 = "HELD Senate - Appropriations"
    19 = "HELD Senate - Appropriations"
    20 = "to Assembly"
    21 = "to Assembly"
    22 = "HELD Senate - Appropriations"
    23 = "to Assembly"
    24 = "HELD Senate - Appropriations"
    25 = "to Assembly"
    26 = "HELD Senate - Appropriations"
    27 = "to Assembly"
    28 = "to Assembly"
    29 = "to Assembly"
    30 = "HELD Senate - Judiciary"
    31 = "to Assembly"
    32 = "to Assembly"
    33 = "HELD Senate - Education"
    34 = "to Assembly"
    35 = "to Assembly"
    36 = "Assembly - Insurance"
    37 = "to Assembly"
    38 = "to Assembly"
    39 = "to Assembly"
    40 = "to Assembly"
    41 = "to Assembly"
    42 = "to Assembly"
    43 = "HELD Senate - Appropriations"
    44 = "HELD Senate - Environmental Quality"
    45 = "to Assembly"
    46 = "Assembly - Military and Veterans Affairs"
    47 = "HELD Senate - Appropriations"
    48 = "HELD Senate - Appropriations"
    49 = "to Assembly"
    50 = "to Assembly"
}

foreach ($row in $statusByRow.Keys) {
    $ws.Cells.Item($row, 6).Value = $statusByRow[$row]
}

# New export timestamp -> new sheet name.
$ws.Name = "Bills-2025-06-04_1735"

# Leave the selection where the export script left the cursor.
$ws.Range("F52").Select() | Out-Null
